$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Expand the table with two new rows ---
$lo = $ws.ListObjects.Item(1)
$row5 = $lo.ListRows.Add()
$row6 = $lo.ListRows.Add()

# --- Row 5: 130. Surrounded Regions ---
$ws.Range("A5").Value() = "130. Surrounded Regions"
$ws.Range("B5").Value() = "Medium"
$ws.Range("C5").Value() = "Graph DFS"
$ws.Range("E5").Value() = "https://leetcode.com/problems/surrounded-regions/solutions/691675/c-beginner-friendly-boundary-dfs-inplace/"
$ws.Range("D5").Value() = 'Use "Boundary DFS". Note that DFS is still O(mxn) time from a for loop, as long as we track visited.'

# --- Row 6: 42. Trapping Rain Water ---
$ws.Range("A6").Value() = "42. Trapping Rain Water"
$ws.Range("B6").Value() = "Hard"
$ws.Range("C6").Value() = "Two Pointers"
$ws.Range("D6").Value() = "Cumulative approach. Only consider l_max and r_max at each step, the smaller of the two is the limiting factor. We track water at each step, subtacting the elevation."
$ws.Range("E6").Value() = "https://leetcode.com/problems/trapping-rain-water/solutions/409175/java-detailed-explanations-illustrations-divide-and-conquer-dp-two-pointers/"

# --- Re-apply the "Medium" orange highlight fill that ListRows.Add loses ---
$ws.Range("B5").Interior.Color = 49407

# --- Hyperlinks for the two new rows, added in this order to match the
#     authored workbook (must happen before the red fill below so the
#     Hyperlink style/font is registered ahead of the new red fill style) ---
$h5 = $ws.Hyperlinks.Add($ws.Range("E5"), "https://leetcode.com/problems/surrounded-regions/solutions/691675/c-beginner-friendly-boundary-dfs-inplace/")
$h6 = $ws.Hyperlinks.Add($ws.Range("E6"), "https://leetcode.com/problems/trapping-rain-water/solutions/409175/java-detailed-explanations-illustrations-divide-and-conquer-dp-two-pointers/")

# --- "Hard" cell gets a red fill ---
$ws.Range("B6").Interior.Color = 255

# --- Fix trailing-space link text on the pre-existing rows (3, then 2) ---
$ws.Range("E3").Value() = "https://leetcode.com/problems/subsets/solutions/27281/a-general-approach-to-backtracking-questions-in-java-subsets-permutations-combination-sum-palindrome-partitioning/ "
$h3 = $ws.Hyperlinks.Add($ws.Range("E3"), "https://leetcode.com/problems/subsets/solutions/27281/a-general-approach-to-backtracking-questions-in-java-subsets-permutations-combination-sum-palindrome-partitioning/ ")

$ws.Range("E2").Value() = "https://leetcode.com/problems/koko-eating-bananas/description/ "
$h2 = $ws.Hyperlinks.Add($ws.Range("E2"), "https://leetcode.com/problems/koko-eating-bananas/description/ ")

# --- Window view tweaks (scroll so column E is visible, then select E11) ---
$win = $excel.ActiveWindow
$win.ScrollColumn = 5
$win.Left = 33255
$win.Top = 2595
$win.Width = 21600
$win.Height = 11010
$ws.Range("E11").Select()
